# Buggy Component List - GITBOOK-240 change request
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "M3-12mm Bolt" quantity correction: 44 -> 42
$ws.Range("C32").Value = 42

# Remove the "M3 nylon nut" row entirely (row 43) - this shifts all rows below up by one
# and shrinks the CustomerList table / sheet dimension accordingly.
$ws.Rows("43:43").Delete()

# Clear the now-orphaned "M3-10mm nylon Screw" entry (row 42), leaving the row blank
# (its "any store" note in column E is kept).
$ws.Range("B42").Value = ""
$ws.Range("C42").Value = ""

# Update the view to reflect where the author left the cursor after editing.
$win = $excel.ActiveWindow()
[void]$win.ScrollRow(24)
[void]$win.ScrollColumn(1)
$ws.Range("C40").Select()
